$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New game rows to append at the bottom of the data table (rows 672-681).
# Columns: A=Away team, B=Away Pts, C=Home team, D=Home Pts,
#          E=Overtime, F=Attend., G=Arena, H=Win, I=Loss
$rows = @(
    @("Washington Wizards", 118, "Detroit Pistons", 104, "No", 17832, "Little Caesars Arena", "Washington Wizards", "Detroit Pistons"),
    @("Miami Heat", 109, "New York Knicks", 125, "No", 17832, "Madison Square Garden (IV)", "New York Knicks", "Miami Heat"),
    @("Philadelphia 76ers", 105, "Denver Nuggets", 111, "No", 17832, "Ball Arena", "Denver Nuggets", "Philadelphia 76ers"),
    @("Houston Rockets", 104, "Brooklyn Nets", 106, "No", 17832, "Barclays Center", "Brooklyn Nets", "Houston Rockets"),
    @("Los Angeles Clippers", 115, "Boston Celtics", 96, "No", 17832, "TD Garden", "Los Angeles Clippers", "Boston Celtics"),
    @("Utah Jazz", 134, "Charlotte Hornets", 122, "No", 17832, "Spectrum Center", "Utah Jazz", "Charlotte Hornets"),
    @("New Orleans Pelicans", 117, "Milwaukee Bucks", 141, "No", 17832, "Fiserv Forum", "Milwaukee Bucks", "New Orleans Pelicans"),
    @("Los Angeles Lakers", 145, "Golden State Warriors", 144, "2OT", 17832, "Chase Center", "Los Angeles Lakers", "Golden State Warriors"),
    @("Minnesota Timberwolves", 112, "San Antonio Spurs", 113, "No", 17832, "Frost Bank Center", "San Antonio Spurs", "Minnesota Timberwolves"),
    @("Sacramento Kings", 120, "Dallas Mavericks", 115, "No", 17832, "American Airlines Center", "Sacramento Kings", "Dallas Mavericks")
)

$startRow = 672
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]

    # Match existing score-column number formatting (thousands separator style)
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
}

# Scroll the view so the newly added rows are visible, matching the saved view state.
$ws.Activate()
$topLeft = $ws.Cells.Item(649, 1)
$excel.ActiveWindow.ScrollRow = $topLeft.Row
$excel.ActiveWindow.ScrollColumn = $topLeft.Column
$ws.Range("A681").Select()
